# Update column F ("dSF") values for the specified rows.
# Mapping of row number -> new value, per commit "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    7  = -1
    9  = 3
    12 = 3
    13 = -1
    14 = 5
    15 = -1
    17 = 0
    18 = -1
    22 = 0
    25 = -2
    31 = -2
    34 = 2
    39 = -1
    40 = -1
    42 = 1
    57 = 8
    60 = 4
    61 = -2
    63 = -1
    64 = -5
    66 = 2
    71 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
